$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update column widths (OOXML width = ColumnWidth + 5/6 for Calibri 11 default font)
$ws.Columns.Item(3).ColumnWidth = 7.166666666666667
$ws.Columns.Item(4).ColumnWidth = 6.166666666666667
$ws.Columns.Item(5).ColumnWidth = 7.166666666666667
$ws.Columns.Item(6).ColumnWidth = 7.166666666666667
$ws.Columns.Item(7).ColumnWidth = 7.166666666666667
$ws.Columns.Item(8).ColumnWidth = 7.166666666666667
$ws.Columns.Item(9).ColumnWidth = 7.166666666666667
$ws.Columns.Item(10).ColumnWidth = 7.166666666666667
$ws.Columns.Item(11).ColumnWidth = 7.166666666666667
$ws.Columns.Item(12).ColumnWidth = 7.166666666666667
$ws.Columns.Item(13).ColumnWidth = 7.166666666666667
$ws.Columns.Item(15).ColumnWidth = 7.166666666666667
$ws.Columns.Item(16).ColumnWidth = 7.166666666666667
$ws.Columns.Item(18).ColumnWidth = 6.166666666666667
$ws.Columns.Item(20).ColumnWidth = 8.166666666666666
$ws.Columns.Item(21).ColumnWidth = 7.166666666666667
$ws.Columns.Item(22).ColumnWidth = 7.166666666666667
$ws.Columns.Item(23).ColumnWidth = 7.166666666666667
$ws.Columns.Item(24).ColumnWidth = 7.166666666666667
$ws.Columns.Item(26).ColumnWidth = 7.166666666666667
$ws.Columns.Item(28).ColumnWidth = 7.166666666666667
$ws.Columns.Item(30).ColumnWidth = 7.166666666666667
$ws.Columns.Item(32).ColumnWidth = 7.166666666666667
$ws.Columns.Item(33).ColumnWidth = 6.166666666666667
$ws.Columns.Item(34).ColumnWidth = 7.166666666666667

# Update data cells for rows 2-5 with new values
$ws.Cells.Item(2, 1).Value = 45108.50694444445
$ws.Cells.Item(2, 2).Value = 8.821999999999999
$ws.Cells.Item(2, 3).Value = 6.243
$ws.Cells.Item(2, 4).Value = 2.794
$ws.Cells.Item(2, 5).Value = 19.605
$ws.Cells.Item(2, 6).Value = 14.205
$ws.Cells.Item(2, 7).Value = 5.893
$ws.Cells.Item(2, 8).Value = 19.279
$ws.Cells.Item(2, 9).Value = 10.957
$ws.Cells.Item(2, 10).Value = 4.563
$ws.Cells.Item(2, 11).Value = 5.697
$ws.Cells.Item(2, 12).Value = 7.759
$ws.Cells.Item(2, 13).Value = 8.683999999999999
$ws.Cells.Item(2, 14).Value = 2.881
$ws.Cells.Item(2, 15).Value = 7.143
$ws.Cells.Item(2, 16).Value = 9.409000000000001
$ws.Cells.Item(2, 17).Value = 6.779
$ws.Cells.Item(2, 18).Value = 1.556
$ws.Cells.Item(2, 19).Value = 0.752
$ws.Cells.Item(2, 20).Value = 101.319
$ws.Cells.Item(2, 21).Value = 19.6
$ws.Cells.Item(2, 22).Value = 6.594
$ws.Cells.Item(2, 23).Value = 12.192
$ws.Cells.Item(2, 24).Value = 6.907
$ws.Cells.Item(2, 25).Value = 0.798
$ws.Cells.Item(2, 26).Value = 11.548
$ws.Cells.Item(2, 27).Value = 5.824
$ws.Cells.Item(2, 28).Value = 5.505
$ws.Cells.Item(2, 29).Value = 6.107
$ws.Cells.Item(2, 30).Value = 8.566000000000001
$ws.Cells.Item(2, 31).Value = 2.095
$ws.Cells.Item(2, 32).Value = 16.882
$ws.Cells.Item(2, 33).Value = 3.332
$ws.Cells.Item(2, 34).Value = 8.242000000000001

$ws.Cells.Item(3, 1).Value = 45108.51388888889
$ws.Cells.Item(3, 2).Value = 4.595
$ws.Cells.Item(3, 3).Value = 3.285
$ws.Cells.Item(3, 4).Value = 1.242
$ws.Cells.Item(3, 5).Value = 10.437
$ws.Cells.Item(3, 6).Value = 7.589
$ws.Cells.Item(3, 7).Value = 3.018
$ws.Cells.Item(3, 8).Value = 15.191
$ws.Cells.Item(3, 9).Value = 5.733
$ws.Cells.Item(3, 10).Value = 2.562
$ws.Cells.Item(3, 11).Value = 2.933
$ws.Cells.Item(3, 12).Value = 4.131
$ws.Cells.Item(3, 13).Value = 4.622
$ws.Cells.Item(3, 14).Value = 1.526
$ws.Cells.Item(3, 15).Value = 3.76
$ws.Cells.Item(3, 16).Value = 5.018
$ws.Cells.Item(3, 17).Value = 3.671
$ws.Cells.Item(3, 18).Value = 0.713
$ws.Cells.Item(3, 19).Value = 0.423
$ws.Cells.Item(3, 20).Value = 49.914
$ws.Cells.Item(3, 21).Value = 10.634
$ws.Cells.Item(3, 22).Value = 3.47
$ws.Cells.Item(3, 23).Value = 6.648
$ws.Cells.Item(3, 24).Value = 3.61
$ws.Cells.Item(3, 25).Value = 0.445
$ws.Cells.Item(3, 26).Value = 8.475
$ws.Cells.Item(3, 27).Value = 3.065
$ws.Cells.Item(3, 28).Value = 2.928
$ws.Cells.Item(3, 29).Value = 3.332
$ws.Cells.Item(3, 30).Value = 4.489
$ws.Cells.Item(3, 31).Value = 0.973
$ws.Cells.Item(3, 32).Value = 14.271
$ws.Cells.Item(3, 33).Value = 1.719
$ws.Cells.Item(3, 34).Value = 4.34

$ws.Cells.Item(4, 1).Value = 45108.52083333334
$ws.Cells.Item(4, 2).Value = 21.94
$ws.Cells.Item(4, 3).Value = 16.384
$ws.Cells.Item(4, 4).Value = 1.437
$ws.Cells.Item(4, 5).Value = 48.077
$ws.Cells.Item(4, 6).Value = 39.06
$ws.Cells.Item(4, 7).Value = 16.845
$ws.Cells.Item(4, 8).Value = 62.578
$ws.Cells.Item(4, 9).Value = 26.685
$ws.Cells.Item(4, 10).Value = 12.021
$ws.Cells.Item(4, 11).Value = 17.271
$ws.Cells.Item(4, 12).Value = 19.262
$ws.Cells.Item(4, 13).Value = 20.555
$ws.Cells.Item(4, 14).Value = 5.767
$ws.Cells.Item(4, 15).Value = 17.294
$ws.Cells.Item(4, 16).Value = 24.442
$ws.Cells.Item(4, 17).Value = 14.72
$ws.Cells.Item(4, 18).Value = 0.626
$ws.Cells.Item(4, 19).Value = 0.851
$ws.Cells.Item(4, 20).Value = 255.909
$ws.Cells.Item(4, 21).Value = 48.261
$ws.Cells.Item(4, 22).Value = 15.963
$ws.Cells.Item(4, 23).Value = 32.302
$ws.Cells.Item(4, 24).Value = 17.061
$ws.Cells.Item(4, 25).Value = 2.244
$ws.Cells.Item(4, 26).Value = 31.968
$ws.Cells.Item(4, 27).Value = 14.1
$ws.Cells.Item(4, 28).Value = 12.539
$ws.Cells.Item(4, 29).Value = 14.7
$ws.Cells.Item(4, 30).Value = 20.328
$ws.Cells.Item(4, 31).Value = 0.644
$ws.Cells.Item(4, 32).Value = 56.903
$ws.Cells.Item(4, 33).Value = 8.917999999999999
$ws.Cells.Item(4, 34).Value = 19.958

$ws.Cells.Item(5, 1).Value = 45108.52777777778
$ws.Cells.Item(5, 2).Value = 23.41
$ws.Cells.Item(5, 3).Value = 17.52
$ws.Cells.Item(5, 4).Value = 1.3
$ws.Cells.Item(5, 5).Value = 51.21
$ws.Cells.Item(5, 6).Value = 41.83
$ws.Cells.Item(5, 7).Value = 18.1
$ws.Cells.Item(5, 8).Value = 71.88
$ws.Cells.Item(5, 9).Value = 28.44
$ws.Cells.Item(5, 10).Value = 12.85
$ws.Cells.Item(5, 11).Value = 18.59
$ws.Cells.Item(5, 12).Value = 20.53
$ws.Cells.Item(5, 13).Value = 21.86
$ws.Cells.Item(5, 14).Value = 6.08
$ws.Cells.Item(5, 15).Value = 18.42
$ws.Cells.Item(5, 16).Value = 26.14
$ws.Cells.Item(5, 17).Value = 15.56
$ws.Cells.Item(5, 18).Value = 0.53
$ws.Cells.Item(5, 19).Value = 0.85
$ws.Cells.Item(5, 20).Value = 273.1
$ws.Cells.Item(5, 21).Value = 51.53
$ws.Cells.Item(5, 22).Value = 17
$ws.Cells.Item(5, 23).Value = 34.61
$ws.Cells.Item(5, 24).Value = 18.21
$ws.Cells.Item(5, 25).Value = 2.4
$ws.Cells.Item(5, 26).Value = 35.59
$ws.Cells.Item(5, 27).Value = 15.02
$ws.Cells.Item(5, 28).Value = 13.31
$ws.Cells.Item(5, 29).Value = 15.63
$ws.Cells.Item(5, 30).Value = 21.63
$ws.Cells.Item(5, 31).Value = 0.48
$ws.Cells.Item(5, 32).Value = 65.47
$ws.Cells.Item(5, 33).Value = 9.550000000000001
$ws.Cells.Item(5, 34).Value = 21.26

# Remove old row 6 (data no longer present in target)
$ws.Rows.Item(6).Delete()

Write-Output "DONE"